# Weekly update: insert two new price records (row 818 and 819) for the
# "Feria Lagunitas de Puerto Montt - Plátano" consolidated sheet.
# Existing rows 818..886 shift down to 820..888.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 818 (pushes everything from 818 down by 2).
$ws.Rows.Item(818).Insert()
$ws.Rows.Item(818).Insert()

# --- New row 818 ---
$ws.Range("A818").Value = 4
$ws.Range("B818").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C818").Value = "Los Lagos"
$ws.Range("D818").Value = 45106
$ws.Range("E818").Value = 10
$ws.Range("F818").Value = "Fruta"
$ws.Range("G818").Value = 100108
$ws.Range("H818").Value = "Tropicales y subtropicales"
$ws.Range("I818").Value = 100108006
$ws.Range("J818").Value = "Plátano"
$ws.Range("K818").Value = "Sin especificar"
$ws.Range("L818").Value = "Pintón"
$ws.Range("M818").Value = 300
$ws.Range("N818").Value = 17000
$ws.Range("O818").Value = 17000
$ws.Range("P818").Value = 17000
$ws.Range("Q818").Value = '$/caja 20 kilos'
$ws.Range("R818").Value = "Ecuador"
$ws.Range("S818").Value = 850
$ws.Range("T818").Value = 20

# --- New row 819 ---
$ws.Range("A819").Value = 4
$ws.Range("B819").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C819").Value = "Los Lagos"
$ws.Range("D819").Value = 45106
$ws.Range("E819").Value = 10
$ws.Range("F819").Value = "Fruta"
$ws.Range("G819").Value = 100108
$ws.Range("H819").Value = "Tropicales y subtropicales"
$ws.Range("I819").Value = 100108006
$ws.Range("J819").Value = "Plátano"
$ws.Range("K819").Value = "Sin especificar"
$ws.Range("L819").Value = "Primera Pintón"
$ws.Range("M819").Value = 800
$ws.Range("N819").Value = 18000
$ws.Range("O819").Value = 19000
$ws.Range("P819").Value = 18500
$ws.Range("Q819").Value = '$/caja 20 kilos'
$ws.Range("R819").Value = "Ecuador"
$ws.Range("S819").Value = 925
$ws.Range("T819").Value = 20
